$d = $word.ActiveDocument

# --- Locate the two paragraphs involved in the edit -----------------------
# "Launch pgAdmin4" is the numbered list item right before the blank
# paragraphs at the tail of the document. The paragraph that needs to be
# removed is the *second* of the two consecutive blank paragraphs
# (ind=360, no list formatting) that sit right before it.
$launchIndex = 0
$blankIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Launch pgAdmin4 ") {
        $launchIndex = $i
        break
    }
}
$blankIndex = $launchIndex - 1

# --- Remove the stray blank paragraph just before "Launch pgAdmin4" -------
$d.Paragraphs.Item($blankIndex).Range.Delete()

# After the delete, "Launch pgAdmin4" shifted up by one.
$launchIndex = $launchIndex - 1
$pLaunch = $d.Paragraphs.Item($launchIndex)

# --- Add a new blank "List Paragraph" line right after it -----------------
# Insert the new paragraph's XML directly (positioned just before the
# paragraph mark of "Launch pgAdmin4") so it picks up a clean
# ListParagraph style with no inherited numbering (w:numPr).
$insertPos = $pLaunch.Range.End - 1
$insertionPoint = $d.Range($insertPos, $insertPos)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p></w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newParaXml)
